$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "B3"   = 6.428999999999999
    "E3"   = 16.337
    "A12"  = -21.533
    "B14"  = 5.944
    "B26"  = 6.516000000000001
    "E30"  = 16.135
    "B31"  = 6.287
    "A32"  = -21.347
    "B35"  = 8.095000000000001
    "A36"  = -20.945
    "B37"  = 8.260000000000002
    "A38"  = -20.354
    "E44"  = 16.16
    "B45"  = 5.786
    "A46"  = -21.497
    "A54"  = -22.209
    "A55"  = -22.164
    "B57"  = 5.331999999999999
    "E58"  = 16.539
    "A67"  = -21.6
    "A69"  = -21.651
    "A72"  = -21.55
    "E84"  = 16.556
    "E89"  = 17.149
    "A91"  = -21.508
    "E91"  = 17.39
    "E92"  = 17.381
    "A99"  = -20.828
    "B100" = 5.558
    "B102" = 7.774000000000001
    "E102" = 16.42
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
